$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-28 Monday" "2025-04-29 Tuesday"

Replace-Text "993÷2=496, 1" "175÷6=29, 1"
Replace-Text "296÷7=42, 2" "567÷9=63, 0"
Replace-Text "310÷8=38, 6" "588÷2=294, 0"
Replace-Text "863÷6=143, 5" "853÷2=426, 1"
Replace-Text "390÷8=48, 6" "349÷8=43, 5"
Replace-Text "805÷3=268, 1" "499÷7=71, 2"
Replace-Text "295÷6=49, 1" "283÷7=40, 3"
Replace-Text "176÷3=58, 2" "864÷5=172, 4"
Replace-Text "672÷3=224, 0" "387÷2=193, 1"
Replace-Text "427÷5=85, 2" "355÷7=50, 5"
Replace-Text "725÷6=120, 5" "599÷7=85, 4"
Replace-Text "519÷9=57, 6" "970÷2=485, 0"
Replace-Text "789÷9=87, 6" "385÷9=42, 7"
Replace-Text "737÷5=147, 2" "343÷7=49, 0"
Replace-Text "550÷7=78, 4" "496÷5=99, 1"
Replace-Text "948÷6=158, 0" "591÷5=118, 1"
Replace-Text "920÷9=102, 2" "289÷7=41, 2"
Replace-Text "930÷6=155, 0" "914÷8=114, 2"
Replace-Text "576÷2=288, 0" "509÷4=127, 1"
Replace-Text "524÷8=65, 4" "302÷3=100, 2"
Replace-Text "684÷4=171, 0" "192÷6=32, 0"
Replace-Text "110÷6=18, 2" "409÷7=58, 3"
Replace-Text "854÷4=213, 2" "644÷2=322, 0"
Replace-Text "451÷3=150, 1" "165÷3=55, 0"
Replace-Text "158÷7=22, 4" "143÷3=47, 2"

Write-Output "Done"
